# Weekly data refresh: insert a new week's row for Cilantro (Terminal La
# Palmera de La Serena) at row 68, pushing the existing rows 68-162 down
# to 69-163. The new row reuses the static descriptive columns of the
# (old) row 68 and its K/L/M/P price-range figures, only the date (D) and
# volume (J) are genuinely new values for the added week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 68..162 down by inserting a blank row at position 68.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new week's record.
$ws.Cells.Item(68, 1).Value = 8
$ws.Cells.Item(68, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(68, 3).Value = 'Coquimbo'
$ws.Cells.Item(68, 4).Value = 44791
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = 100112040
$ws.Cells.Item(68, 7).Value = 'Cilantro'
$ws.Cells.Item(68, 8).Value = 'Sin especificar'
$ws.Cells.Item(68, 9).Value = 'Primera'
$ws.Cells.Item(68, 10).Value = 2800
$ws.Cells.Item(68, 11).Value = 2000
$ws.Cells.Item(68, 12).Value = 2500
$ws.Cells.Item(68, 13).Value = 2250
$ws.Cells.Item(68, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(68, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(68, 16).Value = 1500
$ws.Cells.Item(68, 17).Value = 1.5
$ws.Cells.Item(68, 18).Value = 'Hortaliza'
